$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the data of the observation rows 2-8: the content that
# used to live in row N now lives two rows earlier (wrapping within the
# 2..8 block), i.e. new row R ends up holding the data that used to be in
# row R+2 (mod 7, staying inside 2..8). Row 1 (headers) and row 9 are
# untouched.
#
# Strategy: stage a verbatim copy of rows 2-8 far away (rows 1000-1006) so
# the originals survive being overwritten, then paste each staged row back
# into its new slot. Copying whole cells (rather than re-typing values)
# keeps numbers/text/booleans exactly as they were (e.g. "Antal" values
# stored as text stay text instead of becoming numbers).
#
# PasteSpecial leaves a destination cell alone when the corresponding
# source cell is blank/absent (it never blanks out pre-existing content),
# so a handful of destination cells whose old value has no counterpart in
# the incoming row need to be cleared explicitly first.

# 1) Stage original rows 2-8 at rows 1000-1006 (row 1000 <- row 2, ... row 1006 <- row 8)
$ws.Range("A2:AY8").Copy()
$ws.Range("A1000:AY1006").PasteSpecial(-4163)

# 2) Clear the one cell per destination row whose old value has no
#    corresponding (non-blank) source cell in the incoming data, so it
#    doesn't linger (PasteSpecial leaves a destination cell untouched when
#    the matching source cell is blank/absent).
$ws.Range("M2").Clear()
$ws.Range("M3").Clear()
$ws.Range("I4").Clear()
$ws.Range("I5").Clear()
$ws.Range("K6").Clear()
$ws.Range("K7").Clear()

# 3) Map: new row -> staged row holding the data that should land there
#    new 2 <- old 4 (staged 1002)
#    new 3 <- old 5 (staged 1003)
#    new 4 <- old 6 (staged 1004)
#    new 5 <- old 7 (staged 1005)
#    new 6 <- old 8 (staged 1006)
#    new 7 <- old 2 (staged 1000)
#    new 8 <- old 3 (staged 1001)
$ws.Range("A1002:AY1002").Copy()
$ws.Range("A2:AY2").PasteSpecial(-4163)

$ws.Range("A1003:AY1003").Copy()
$ws.Range("A3:AY3").PasteSpecial(-4163)

$ws.Range("A1004:AY1004").Copy()
$ws.Range("A4:AY4").PasteSpecial(-4163)

$ws.Range("A1005:AY1005").Copy()
$ws.Range("A5:AY5").PasteSpecial(-4163)

$ws.Range("A1006:AY1006").Copy()
$ws.Range("A6:AY6").PasteSpecial(-4163)

$ws.Range("A1000:AY1000").Copy()
$ws.Range("A7:AY7").PasteSpecial(-4163)

$ws.Range("A1001:AY1001").Copy()
$ws.Range("A8:AY8").PasteSpecial(-4163)

# 4) Remove the staging area
$ws.Range("A1000:AY1006").Clear()
